# Update crypto Price (D) and Volume(1h) (E) columns to refreshed values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.845.43'
$ws.Range('E2').Value = '  -0.02%  '
$ws.Range('D3').Value = '2.411.64'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = "'550.91"
$ws.Range('E5').Value = '  -0.68%  '
$ws.Range('D6').Value = "'136.79"
$ws.Range('E6').Value = '  -0.80%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  +4.23%  '
$ws.Range('E9').Value = '  -1.47%  '
$ws.Range('E10').Value = '  -2.28%  '
$ws.Range('E11').Value = '  -0.94%  '
$ws.Range('E12').Value = '  -1.86%  '
$ws.Range('D13').Value = "'25.27"
$ws.Range('E13').Value = '  +2.33%  '
$ws.Range('D14').Value = '2.839.23'
$ws.Range('E14').Value = '  -0.53%  '
$ws.Range('D15').Value = '59.762.24'
$ws.Range('E15').Value = '  -0.02%  '
$ws.Range('E16').Value = '  -1.77%  '
$ws.Range('D17').Value = '2.411.17'
$ws.Range('E17').Value = '  -1.43%  '
$ws.Range('D18').Value = "'11.32"
$ws.Range('E18').Value = '  -0.42%  '
$ws.Range('D19').Value = "'4.41"
$ws.Range('E19').Value = '  -0.53%  '
$ws.Range('D20').Value = "'329.01"
$ws.Range('E20').Value = '  -1.65%  '
$ws.Range('E21').Value = '  -4.06%  '
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').Value = "'66.27"
$ws.Range('E23').Value = '  +2.56%  '
$ws.Range('E24').Value = '  +1.92%  '
$ws.Range('D25').Value = "'8.63"
$ws.Range('E25').Value = '  +0.23%  '
$ws.Range('E27').Value = '  -0.60%  '
$ws.Range('D28').Value = '0.0₃0770'
$ws.Range('E28').Value = '  -2.24%  '
$ws.Range('E29').Value = '  -2.17%  '
$ws.Range('D30').Value = "'168.95"
$ws.Range('E30').Value = '  -1.12%  '
$ws.Range('E31').Value = '  -4.31%  '
$ws.Range('D32').Value = "'18.60"
$ws.Range('E32').Value = '  -0.41%  '
$ws.Range('E33').Value = '  -1.28%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('E35').Value = '  -0.78%  '
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range('D37').Value = "'4.18"
$ws.Range('E37').Value = '  -1.76%  '
$ws.Range('E38').Value = '  -2.14%  '
$ws.Range('D39').Value = "'320.51"
$ws.Range('D40').Value = "'0.405"
$ws.Range('E40').Value = '  -4.01%  '
$ws.Range('D41').Value = "'3.65"
$ws.Range('E41').Value = '  -2.43%  '
$ws.Range('D42').Value = "'140.01"
$ws.Range('E42').Value = '  -1.76%  '
$ws.Range('E43').Value = '  +0.31%  '
$ws.Range('D44').Value = "'19.48"
$ws.Range('E44').Value = '  +1.28%  '
$ws.Range('D45').Value = "'0.0514"
$ws.Range('E45').Value = '  -1.56%  '
$ws.Range('D46').Value = "'0.577"
$ws.Range('E46').Value = '  +1.27%  '
$ws.Range('E47').Value = '  -1.51%  '
$ws.Range('D48').Value = "'0.386"
$ws.Range('E48').Value = '  -6.08%  '
$ws.Range('D49').Value = "'11.04"
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('E50').Value = '  -3.29%  '
$ws.Range('D51').Value = "'4.67"
$ws.Range('E51').Value = '  -1.14%  '
